# Atualização agendada das bases de dados
# Appends 2025 monthly data (MAR and TERRA, SERGIPE) below the existing
# historical series, extending the sheet from row 673 to row 697.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (date text, value) pairs, in month order Jan..Dec 2025
$marValues = @(558.401, 150.578, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$terraValues = @(58991.751, 51664.747, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

$startRow = 674

# MAR block: rows 674-685
for ($i = 0; $i -lt 12; $i++) {
    $r = $startRow + $i
    $month = $i + 1
    $dateText = "'{0:D2}/{1:D2}/2025" -f 1, $month

    # Leading "'" forces the date-looking text to stay a literal string
    # instead of Excel auto-converting it to a date serial number.
    # ClearFormats() then strips the "quote prefix" cell style Excel
    # applies for that, so the cell keeps the default (unstyled) look
    # matching the rest of the column.
    $ws.Cells.Item($r, 1).Value = $dateText
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).Value = "MAR"
    $ws.Cells.Item($r, 3).Value = $marValues[$i]
    $ws.Cells.Item($r, 4).Value = "SERGIPE"
}

# TERRA block: rows 686-697
$terraStartRow = $startRow + 12
for ($i = 0; $i -lt 12; $i++) {
    $r = $terraStartRow + $i
    $month = $i + 1
    $dateText = "'{0:D2}/{1:D2}/2025" -f 1, $month

    $ws.Cells.Item($r, 1).Value = $dateText
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).Value = "TERRA"
    $ws.Cells.Item($r, 3).Value = $terraValues[$i]
    $ws.Cells.Item($r, 4).Value = "SERGIPE"
}
